$wb = $excel.ActiveWorkbook

# --- AddOpportunity sheet: swap the LOB-related test data for the row ---
$ws1 = $wb.Worksheets.Item("AddOpportunity")

# IndustryGroup / Sector pairing changed from
#   "BUS - Business Services" / "Dealership & Rental Services"
# to
#   "HC - Healthcare" / "Dental"
$ws1.Range("D2").Value = "HC - Healthcare"
$ws1.Range("E2").Value = "Dental"

# Sector cell picks up wrap + vertical-center formatting
$ws1.Range("E2").WrapText = $true
$ws1.Range("E2").VerticalAlignment = -4108

# Staff changed from "Drew Koecher" to "Karan Chopra"
$ws1.Range("N2").Value = "Karan Chopra"

# --- Users sheet: same staff-name swap ---
$ws2 = $wb.Worksheets.Item("Users")
$ws2.Range("A2").Value = "Karan Chopra"

# --- restore on-screen selections to match the saved state ---
$ws2.Range("C15").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D9").Select() | Out-Null
